$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date (as an Excel date
# serial) for every data row (rows 2-272). The automatic export bumped this
# stamp forward by a single day (46061 -> 46062, i.e. 2026-02-08 -> 2026-02-09)
# for every row in the sheet.
$ws.Range("C2:C272").Value = 46062
